# Code Refactor as per Sonarlint
$wb = $excel.ActiveWorkbook

$wsRunManager = $wb.Worksheets.Item("RunManager")
$wsTestData = $wb.Worksheets.Item("TestData")

# Remember which sheet is active so we can restore it after touching the
# other sheet's selection (changing a range's selection on a worksheet
# normally activates that worksheet in Excel).
$originalActive = $wb.ActiveSheet.Name

# Update the stored selection on the RunManager sheet (was E3 -> now D16)
$wsRunManager.Activate()
$wsRunManager.Range("D16").Select()

# Restore the original active sheet (TestData) so the active tab doesn't change
$wb.Worksheets.Item($originalActive).Activate()

# Update the two corrected/refactored page-title strings used in TestData
$wsTestData.Range("I5").Value = "Amazon.com: Computers & Accessories: Electronics: Computer Accessories & Peripherals, T3ablet Accessories & More"
$wsTestData.Range("I7").Value = "Amazon.com: Computers & Accessories: Electronics: Co-mputer Accessories & Peripherals, Tablet Accessories & More"
